# Auto-generated edit script: update PaymentType/Date cells (column B) with new
# Katalon-generated timestamps for Personal_SearchTransaction_Generic_TC test data.
$wb = $excel.ActiveWorkbook

# --- Estimated ---
$wsEstimated = $wb.Worksheets.Item("Estimated")
$wsEstimated.Range("B2").Value = "Fri Mar 08 23:24:32 EST 2024"
$wsEstimated.Range("B3").Value = "Fri Mar 08 23:25:23 EST 2024"
$wsEstimated.Range("B4").Value = "Fri Mar 08 23:26:16 EST 2024"
$wsEstimated.Range("B5").Value = "Fri Mar 08 23:27:08 EST 2024"
$wsEstimated.Range("B6").Value = "Fri Mar 08 23:28:01 EST 2024"
$wsEstimated.Range("B7").Value = "Fri Mar 08 23:28:54 EST 2024"

# --- Existing ---
$wsExisting = $wb.Worksheets.Item("Existing")
$wsExisting.Range("B2").Value = "Fri Mar 08 23:29:46 EST 2024"
$wsExisting.Range("B3").Value = "Fri Mar 08 23:30:38 EST 2024"
$wsExisting.Range("B4").Value = "Fri Mar 08 23:31:30 EST 2024"
$wsExisting.Range("B5").Value = "Fri Mar 08 23:32:22 EST 2024"
$wsExisting.Range("B6").Value = "Fri Mar 08 23:33:15 EST 2024"
$wsExisting.Range("B7").Value = "Fri Mar 08 23:34:07 EST 2024"
$wsExisting.Range("B8").Value = "Fri Mar 08 23:35:00 EST 2024"
$wsExisting.Range("B9").Value = "Fri Mar 08 23:35:54 EST 2024"
$wsExisting.Range("B10").Value = "Fri Mar 08 23:36:48 EST 2024"
$wsExisting.Range("B11").Value = "Fri Mar 08 23:37:42 EST 2024"
$wsExisting.Range("B12").Value = "Fri Mar 08 23:38:34 EST 2024"
$wsExisting.Range("B13").Value = "Fri Mar 08 23:39:26 EST 2024"
$wsExisting.Range("B14").Value = "Fri Mar 08 23:40:19 EST 2024"
$wsExisting.Range("B15").Value = "Fri Mar 08 23:41:12 EST 2024"
$wsExisting.Range("B16").Value = "Fri Mar 08 23:42:05 EST 2024"
$wsExisting.Range("B17").Value = "Fri Mar 08 23:42:57 EST 2024"
$wsExisting.Range("B18").Value = "Fri Mar 08 23:43:49 EST 2024"
$wsExisting.Range("B19").Value = "Fri Mar 08 23:44:42 EST 2024"

# --- NewTaxReturn ---
$wsNewTaxReturn = $wb.Worksheets.Item("NewTaxReturn")
$wsNewTaxReturn.Range("B2").Value = "Fri Mar 08 23:45:34 EST 2024"
$wsNewTaxReturn.Range("B3").Value = "Fri Mar 08 23:46:26 EST 2024"
$wsNewTaxReturn.Range("B4").Value = "Fri Mar 08 23:47:17 EST 2024"
$wsNewTaxReturn.Range("B5").Value = "Fri Mar 08 23:48:09 EST 2024"
$wsNewTaxReturn.Range("B6").Value = "Fri Mar 08 23:49:00 EST 2024"
$wsNewTaxReturn.Range("B7").Value = "Fri Mar 08 23:49:52 EST 2024"
$wsNewTaxReturn.Range("B8").Value = "Fri Mar 08 23:50:44 EST 2024"
$wsNewTaxReturn.Range("B9").Value = "Fri Mar 08 23:51:36 EST 2024"
$wsNewTaxReturn.Range("B10").Value = "Fri Mar 08 23:52:28 EST 2024"
$wsNewTaxReturn.Range("B11").Value = "Fri Mar 08 23:53:20 EST 2024"
$wsNewTaxReturn.Range("B12").Value = "Fri Mar 08 23:54:13 EST 2024"
$wsNewTaxReturn.Range("B13").Value = "Fri Mar 08 23:55:05 EST 2024"
$wsNewTaxReturn.Range("B14").Value = "Fri Mar 08 23:55:56 EST 2024"
$wsNewTaxReturn.Range("B15").Value = "Fri Mar 08 23:56:48 EST 2024"
$wsNewTaxReturn.Range("B16").Value = "Fri Mar 08 23:57:41 EST 2024"
$wsNewTaxReturn.Range("B17").Value = "Fri Mar 08 23:58:32 EST 2024"
$wsNewTaxReturn.Range("B18").Value = "Fri Mar 08 23:59:26 EST 2024"
$wsNewTaxReturn.Range("B19").Value = "Sat Mar 09 00:00:20 EST 2024"
$wsNewTaxReturn.Range("B20").Value = "Sat Mar 09 00:01:14 EST 2024"
$wsNewTaxReturn.Range("B21").Value = "Sat Mar 09 00:02:06 EST 2024"
$wsNewTaxReturn.Range("B22").Value = "Sat Mar 09 00:02:59 EST 2024"
$wsNewTaxReturn.Range("B23").Value = "Sat Mar 09 00:03:52 EST 2024"
$wsNewTaxReturn.Range("B24").Value = "Sat Mar 09 00:04:46 EST 2024"
$wsNewTaxReturn.Range("B25").Value = "Sat Mar 09 00:05:40 EST 2024"
$wsNewTaxReturn.Range("B26").Value = "Sat Mar 09 00:06:33 EST 2024"
$wsNewTaxReturn.Range("B27").Value = "Sat Mar 09 00:07:27 EST 2024"
$wsNewTaxReturn.Range("B28").Value = "Sat Mar 09 00:08:20 EST 2024"
$wsNewTaxReturn.Range("B29").Value = "Sat Mar 09 00:09:13 EST 2024"
$wsNewTaxReturn.Range("B30").Value = "Sat Mar 09 00:10:07 EST 2024"
$wsNewTaxReturn.Range("B31").Value = "Sat Mar 09 00:11:00 EST 2024"
$wsNewTaxReturn.Range("B32").Value = "Sat Mar 09 00:11:54 EST 2024"
$wsNewTaxReturn.Range("B33").Value = "Sat Mar 09 00:12:46 EST 2024"
$wsNewTaxReturn.Range("B34").Value = "Sat Mar 09 00:13:39 EST 2024"
$wsNewTaxReturn.Range("B35").Value = "Sat Mar 09 00:14:32 EST 2024"
$wsNewTaxReturn.Range("B36").Value = "Sat Mar 09 00:15:25 EST 2024"
$wsNewTaxReturn.Range("B37").Value = "Sat Mar 09 00:16:19 EST 2024"
$wsNewTaxReturn.Range("B38").Value = "Sat Mar 09 00:17:12 EST 2024"
$wsNewTaxReturn.Range("B39").Value = "Sat Mar 09 00:18:06 EST 2024"
$wsNewTaxReturn.Range("B40").Value = "Sat Mar 09 00:18:59 EST 2024"
$wsNewTaxReturn.Range("B41").Value = "Sat Mar 09 00:19:53 EST 2024"
$wsNewTaxReturn.Range("B42").Value = "Sat Mar 09 00:20:47 EST 2024"
$wsNewTaxReturn.Range("B43").Value = "Sat Mar 09 00:21:40 EST 2024"
$wsNewTaxReturn.Range("B44").Value = "Sat Mar 09 00:22:34 EST 2024"
$wsNewTaxReturn.Range("B45").Value = "Sat Mar 09 00:23:25 EST 2024"
$wsNewTaxReturn.Range("B46").Value = "Sat Mar 09 00:24:20 EST 2024"
$wsNewTaxReturn.Range("B47").Value = "Sat Mar 09 00:25:14 EST 2024"
$wsNewTaxReturn.Range("B48").Value = "Sat Mar 09 00:26:07 EST 2024"
$wsNewTaxReturn.Range("B49").Value = "Sat Mar 09 00:27:01 EST 2024"
$wsNewTaxReturn.Range("B50").Value = "Sat Mar 09 00:27:54 EST 2024"
$wsNewTaxReturn.Range("B51").Value = "Sat Mar 09 00:28:47 EST 2024"
$wsNewTaxReturn.Range("B52").Value = "Sat Mar 09 00:29:41 EST 2024"

# --- Personal_EL ---
$wsPersonalEL = $wb.Worksheets.Item("Personal_EL")
$wsPersonalEL.Range("B2").Value = "Sat Mar 09 00:30:35 EST 2024"
$wsPersonalEL.Range("B3").Value = "Sat Mar 09 00:31:25 EST 2024"

# --- Personal_IND ---
$wsPersonalIND = $wb.Worksheets.Item("Personal_IND")
$wsPersonalIND.Range("B2").Value = "Sat Mar 09 00:32:16 EST 2024"
$wsPersonalIND.Range("B4").Value = "Sat Mar 09 00:33:06 EST 2024"
$wsPersonalIND.Range("B5").Value = "Sat Mar 09 00:33:56 EST 2024"
$wsPersonalIND.Range("B6").Value = "Sat Mar 09 00:34:45 EST 2024"
$wsPersonalIND.Range("B7").Value = "Sat Mar 09 00:35:36 EST 2024"
$wsPersonalIND.Range("B8").Value = "Sat Mar 09 00:36:25 EST 2024"
$wsPersonalIND.Range("B9").Value = "Sat Mar 09 00:37:15 EST 2024"

# --- Personal_JNT ---
$wsPersonalJNT = $wb.Worksheets.Item("Personal_JNT")
$wsPersonalJNT.Range("B2").Value = "Sat Mar 09 00:38:05 EST 2024"
$wsPersonalJNT.Range("B4").Value = "Sat Mar 09 00:39:00 EST 2024"
$wsPersonalJNT.Range("B5").Value = "Sat Mar 09 00:39:55 EST 2024"
$wsPersonalJNT.Range("B6").Value = "Sat Mar 09 00:40:51 EST 2024"

